$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 holds the 0-ohm resistor group (description 301010292).
# Add R1 and R3 to the designator list and bump the quantity 5 -> 7.
$ws.Range("A8").Value = "R1, R3, R8, R12, R15, R17, R41"
$ws.Range("C8").Value = 7

# Update the view state to match: scroll back to top-left and select C8.
$ws.Range("C8").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
